$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing Arpit "out" time and Nic's lunch-out time for 2018-07-03 (row 11)
$ws.Range("C11").Value = 0.20833333333333334
$ws.Range("E11").Value = 0.23611111111111113
$ws.Range("C11").NumberFormat = "h:mm"
$ws.Range("E11").NumberFormat = "h:mm"

# Add row 12 (2018-07-04)
$ws.Range("A12").Value = 43285
$ws.Range("B12").Value = 0.38194444444444442
$ws.Range("C12").Value = 0.20833333333333334
$ws.Range("D12").Value = 0.38541666666666669
$ws.Range("E12").Value = 0.20833333333333334
$ws.Range("B12:E12").NumberFormat = "h:mm"

# Add row 13 (2018-07-05)
$ws.Range("A13").Value = 43286
$ws.Range("B13").Value = 0.38541666666666669
$ws.Range("C13").Value = 0.20833333333333334
$ws.Range("D13").Value = 0.375
$ws.Range("E13").Value = 0.25
$ws.Range("B13:E13").NumberFormat = "h:mm"

# Update the active selection to match the saved view state
$ws.Range("E15").Select()
